$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A7").Value = "32d33a1d-9237-4406-9595-aff90a885b57.md"
$overview.Range("D7").Value = "2016-03-22 08:23:17"
$overview.Range("A8").Value = "49ebac29-deff-46b1-9ea7-0b70b6c7f2ee.md"
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("A9").Value = "a555daca-8025-49ee-a341-49a5b74d6189.md"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"
$overview.Range("D9").Value = "2016-03-22 08:17:17"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A7").Value = "32d33a1d-9237-4406-9595-aff90a885b57.md"
$zhcn.Range("D7").Value = "32d33a1d-9237-4406-9595-aff90a885b57.09aca4bb52a8223fb1188a5a6a266e3861b9067d.zh-cn.xlf"
$zhcn.Range("E7").Value = "2016-03-22 08:23:13"
$zhcn.Range("A8").Value = "49ebac29-deff-46b1-9ea7-0b70b6c7f2ee.md"
$zhcn.Range("C8").Value = "In Translation"
$zhcn.Range("D8").Value = "49ebac29-deff-46b1-9ea7-0b70b6c7f2ee.a27e0f205ea9883823b1aec313924100a3a6b42b.zh-cn.xlf"
$zhcn.Range("A9").Value = "a555daca-8025-49ee-a341-49a5b74d6189.md"
$zhcn.Range("C9").Value = "In Translation"
$zhcn.Range("D9").Value = "a555daca-8025-49ee-a341-49a5b74d6189.9db5307595bc1ae0120a98e94944b584cf8ca52e.zh-cn.xlf"
$zhcn.Range("E9").Value = "2016-03-22 08:16:59"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A7").Value = "32d33a1d-9237-4406-9595-aff90a885b57.md"
$dede.Range("D7").Value = "32d33a1d-9237-4406-9595-aff90a885b57.09aca4bb52a8223fb1188a5a6a266e3861b9067d.de-de.xlf"
$dede.Range("E7").Value = "2016-03-22 08:23:17"
$dede.Range("A8").Value = "49ebac29-deff-46b1-9ea7-0b70b6c7f2ee.md"
$dede.Range("C8").Value = "In Translation"
$dede.Range("D8").Value = "49ebac29-deff-46b1-9ea7-0b70b6c7f2ee.a27e0f205ea9883823b1aec313924100a3a6b42b.de-de.xlf"
$dede.Range("A9").Value = "a555daca-8025-49ee-a341-49a5b74d6189.md"
$dede.Range("C9").Value = "In Translation"
$dede.Range("D9").Value = "a555daca-8025-49ee-a341-49a5b74d6189.9db5307595bc1ae0120a98e94944b584cf8ca52e.de-de.xlf"
$dede.Range("E9").Value = "2016-03-22 08:17:17"
